# "online shopping process test pt 1"
#
# - Adds a new worksheet "testOnlineShoppingProcess" right before
#   "testNavMenuLinks" (i.e. as the second-to-last tab), populated with the
#   url/login/password header row + values used by the other login-style
#   sheets (testLogIn, testAddAddress, ...), and makes it the active sheet
#   with the selection parked one column past the data (D1).
# - The previously-active sheet ("testLogIn") loses its "tabSelected" flag
#   and its lingering selection is reset back to A2.

$wb = $excel.ActiveWorkbook

# testLogIn was the active/selected tab before this edit; clear that out and
# reset its selection to A2 (matches the diff: tabSelected removed, D2 -> A2).
$loginSheet = $wb.Worksheets.Item("testLogIn")
$loginSheet.Activate()
$loginSheet.Range("A2").Select()

# Insert the new sheet directly before the last tab (testNavMenuLinks), so it
# ends up second-to-last, same as the target layout.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($lastSheet)
$newSheet.Name = "testOnlineShoppingProcess"

$newSheet.Range("A1").Value = "url"
$newSheet.Range("B1").Value = "login"
$newSheet.Range("C1").Value = "password"
$newSheet.Range("A2").Value = "https://magento.softwaretestingboard.com/customer/account/login/referer/"
$newSheet.Range("B2").Value = "sistulostu@gufum.com"
$newSheet.Range("C2").Value = "password!123"

# Leave the selection one cell past the last data column, matching the diff.
$newSheet.Range("D1").Select()
